# Appends the newly retrieved job-numbers data point as a new row at the
# bottom of the sheet, mirroring the existing table layout/formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162

# Locate the current last populated row in column A (date column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$newRow = $lastRow + 1

# New data row values.
$data = @(44353.79480840814, 74722, 62994, 3309, 2101, 1474, 19741, 1397, 886, 196)

for ($col = 1; $col -le $data.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $data[$col - 1]
}

# Match the date column's number format/style used by the rest of the table.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat

# The source dataset this sheet mirrors re-emitted the previous row's
# timestamp with slightly different floating point precision on this run.
$ws.Cells.Item($lastRow, 1).Value = 44352.79363961806
